$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new row at row 10 - this pushes the existing instrument rows
# (previously 10:32) down to 11:33 and grows the used range to E33.
$ws.Rows.Item(10).Insert()

# Populate the newly inserted row with the new instrument/task entry.
$ws.Cells.Item(10, 1).Value = "Print current point on screen"
$ws.Cells.Item(10, 2).Value = 1
$ws.Cells.Item(10, 3).Value = 3
$ws.Cells.Item(10, 4).Value = 3
$ws.Cells.Item(10, 5).Formula = "=B10*C10/D10"

# Match the center-aligned formatting used by the other data rows in column B.
$ws.Cells.Item(10, 2).HorizontalAlignment = -4108

# The color-scale conditional formats on columns C, D and E covered rows
# 2:32 - grow them in place to 2:33 so the new row is covered too (column B's
# two conditional formats are untouched, matching the target workbook).
$ws.Range("C2:C32").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("C2:C33"))
$ws.Range("D2:D32").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("D2:D33"))
$ws.Range("E2:E32").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("E2:E33"))

# Move the selection back to the top-left of the data, matching the saved
# workbook's cursor position.
$ws.Range("B1").Select()
